$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 21 revised run times (B25: Part 1, C25: Part 2 -> E25 is the shared SUM formula)
$ws.Range("B25").Value = 7.78277159994468
$ws.Range("C25").Value = 0.00063530006445944298

# Day 22 new run times (row 26 was previously empty)
$ws.Range("B26").Value = 27.250088499975298
$ws.Range("C26").Value = 0.0010418000165373

# Reflect the author's final selection/active cell in the saved view state
$ws.Range("I29").Select()
